$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that wrapped "First Finger" in the
#    "1 to 1000" bullet (<w:bookmarkStart .../> ... <w:bookmarkEnd .../>).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
if ($bm -ne $null) {
    $bm.Delete()
}

# ---------------------------------------------------------------------------
# 2. Find the (non-bold) "Define the Problem" paragraph that introduces the
#    "Predicting Fingers" problem-solving checklist, and make it bold.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Define the Problem`r" -and $p.Range.Font.Bold -eq 0) {
        $target = $p
    }
}

$target.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# 3. Insert the explanation paragraph + the (empty, bold) bookmark paragraph
#    right after it, before "Break the problem apart".
# ---------------------------------------------------------------------------
$apos = [char]0x2019

$bodyXmlOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$bodyXmlClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPrFont = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

$explanationText1 = 'This wasn' + $apos + 't as difficult of a problem to solve as the previous two problems encountered in this activity.  It was more a matter of learning a pattern.  The situation involved a little girl counting on her fingers from thumb to pinky and then reversing.  The problem was'
$explanationText2 = ' what finger would she end up on when she got to 10, 100 and 1000.'

$explanationP = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/>' + $rPrFont + '</w:pPr>' + `
    '<w:r>' + $rPrFont + '<w:t>' + $explanationText1 + '</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPrFont + '<w:t>,</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPrFont + '<w:t xml:space="preserve">' + $explanationText2 + '</w:t></w:r>' + `
    '</w:p>'

$bookmarkP = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# Create an empty paragraph right after "Define the Problem" and fill it with
# the explanation text via InsertXML (keeps full control over run/pPr XML).
$target.Range.InsertParagraphAfter()
$explPara = $target.Next()
$explPara.Range.InsertXML($bodyXmlOpen + $explanationP + $bodyXmlClose)

# Create a second empty paragraph after the explanation paragraph and fill it
# with the (empty) bold bookmark paragraph.
$explPara.Range.InsertParagraphAfter()
$bmPara = $explPara.Next()
$bmPara.Range.InsertXML($bodyXmlOpen + $bookmarkP + $bodyXmlClose)

Write-Output "done"
